$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 24531
$ws.Range("E2").Value = 1642
$ws.Range("F2").Value = 1642
$ws.Range("G2").Value = 1224
$ws.Range("H2").Value = 765
$ws.Range("I2").Value = 717
$ws.Range("J2").Value = 48
$ws.Range("K2").Value = 22850
$ws.Range("L2").Value = 11103
$ws.Range("M2").Value = 11747
$ws.Range("N2").Value = 11452
$ws.Range("O2").Value = 295
$ws.Range("P2").Value = 300
$ws.Range("Q2").Value = 83
$ws.Range("R2").Value = -943
$ws.Range("S2").Value = 831
$ws.Range("T2").Value = 292
$ws.Range("U2").Value = -209
$ws.Range("V2").Value = 7178
$ws.Range("W2").Value = 6.69
$ws.Range("X2").Value = 3.12
$ws.Range("Y2").Value = 6.46
$ws.Range("Z2").Value = 3.58
$ws.Range("AA2").Value = 94.52
$ws.Range("AB2").Value = 3592.4
$ws.Range("AC2").Value = 11952
$ws.Range("AD2").Value = 7.39
$ws.Range("AE2").Value = 196971
$ws.Range("AF2").Value = 0.45
$ws.Range("AG2").Value = 1500
$ws.Range("AH2").Value = 1.7
$ws.Range("AI2").Value = 12.16
$ws.Range("AJ2").Value = 6000000

# Row 3
$ws.Range("D3").Value = 21917
$ws.Range("E3").Value = 777
$ws.Range("F3").Value = 777
$ws.Range("G3").Value = 685
$ws.Range("H3").Value = 458
$ws.Range("I3").Value = 457
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 20739
$ws.Range("L3").Value = 8467
$ws.Range("M3").Value = 12271
$ws.Range("N3").Value = 11958
$ws.Range("O3").Value = 314
$ws.Range("P3").Value = 300
$ws.Range("Q3").Value = 2930
$ws.Range("R3").Value = -998
$ws.Range("S3").Value = -2396
$ws.Range("T3").Value = 315
$ws.Range("U3").Value = 2614
$ws.Range("V3").Value = 4945
$ws.Range("W3").Value = 3.55
$ws.Range("X3").Value = 2.09
$ws.Range("Y3").Value = 3.91
$ws.Range("Z3").Value = 2.1
$ws.Range("AA3").Value = 69
$ws.Range("AB3").Value = 3711.49
$ws.Range("AC3").Value = 7618
$ws.Range("AD3").Value = 7.84
$ws.Range("AE3").Value = 205671
$ws.Range("AF3").Value = 0.29
$ws.Range("AG3").Value = 1750
$ws.Range("AH3").Value = 2.93
$ws.Range("AI3").Value = 22.26
$ws.Range("AJ3").Value = 6000000

# Row 4
$ws.Range("D4").Value = 17975
$ws.Range("E4").Value = 772
$ws.Range("F4").Value = 772
$ws.Range("G4").Value = 722
$ws.Range("H4").Value = 653
$ws.Range("I4").Value = 651
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 21697
$ws.Range("L4").Value = 8776
$ws.Range("M4").Value = 12921
$ws.Range("N4").Value = 12596
$ws.Range("O4").Value = 325
$ws.Range("P4").Value = 300
$ws.Range("Q4").Value = 1196
$ws.Range("R4").Value = -1069
$ws.Range("S4").Value = -123
$ws.Range("T4").Value = 1211
$ws.Range("U4").Value = -15
$ws.Range("V4").Value = 4906
$ws.Range("W4").Value = 4.3
$ws.Range("X4").Value = 3.63
$ws.Range("Y4").Value = 5.3
$ws.Range("Z4").Value = 3.08
$ws.Range("AA4").Value = 67.92
$ws.Range("AB4").Value = 3889.64
$ws.Range("AC4").Value = 10854
$ws.Range("AD4").Value = 9.49
$ws.Range("AE4").Value = 216644
$ws.Range("AF4").Value = 0.48
$ws.Range("AG4").Value = 1750
$ws.Range("AH4").Value = 1.7
$ws.Range("AI4").Value = 15.62
$ws.Range("AJ4").Value = 6000000

# Row 5
$ws.Range("D5").Value = 15070
$ws.Range("E5").Value = 522
$ws.Range("F5").Value = 522
$ws.Range("G5").Value = 201
$ws.Range("H5").Value = 291
$ws.Range("I5").Value = 255
$ws.Range("J5").Value = 36
$ws.Range("K5").Value = 24117
$ws.Range("L5").Value = 11426
$ws.Range("M5").Value = 12691
$ws.Range("N5").Value = 12370
$ws.Range("O5").Value = 321
$ws.Range("P5").Value = 300
$ws.Range("Q5").Value = -1415
$ws.Range("R5").Value = -751
$ws.Range("S5").Value = 2229
$ws.Range("T5").Value = 654
$ws.Range("U5").Value = -2069
$ws.Range("V5").Value = 6974
$ws.Range("W5").Value = 3.47
$ws.Range("X5").Value = 1.93
$ws.Range("Y5").Value = 2.04
$ws.Range("Z5").Value = 1.27
$ws.Range("AA5").Value = 90.03
$ws.Range("AB5").Value = 3943.81
$ws.Range("AC5").Value = 4248
$ws.Range("AD5").Value = 25.19
$ws.Range("AE5").Value = 212768
$ws.Range("AF5").Value = 0.5
$ws.Range("AG5").Value = 1800
$ws.Range("AH5").Value = 1.68
$ws.Range("AI5").Value = 41.06
$ws.Range("AJ5").Value = 6000000

# Row 6
$ws.Range("D6").Value = 17817
$ws.Range("E6").Value = 805
$ws.Range("F6").Value = 805
$ws.Range("G6").Value = 1924
$ws.Range("H6").Value = 3228
$ws.Range("I6").Value = 3164
$ws.Range("K6").Value = 24454
$ws.Range("L6").Value = 10928
$ws.Range("M6").Value = 13526
$ws.Range("N6").Value = 9413
$ws.Range("P6").Value = 207
$ws.Range("Q6").Value = 1275
$ws.Range("R6").Value = -481
$ws.Range("S6").Value = -564
$ws.Range("T6").Value = 496
$ws.Range("U6").Value = 779
$ws.Range("V6").Value = 6772
$ws.Range("W6").Value = 4.52
$ws.Range("X6").Value = 18.12
$ws.Range("Y6").Value = 29.05
$ws.Range("Z6").Value = 13.29
$ws.Range("AA6").Value = 80.79
$ws.Range("AB6").Value = 4367.77
$ws.Range("AC6").Value = 62315
$ws.Range("AD6").Value = 0.79
$ws.Range("AE6").Value = 232782
$ws.Range("AF6").Value = 0.21
$ws.Range("AG6").Value = 1800
$ws.Range("AH6").Value = 3.65
$ws.Range("AI6").Value = 2.3
$ws.Range("AJ6").Value = 4141657

# Row 7 - clear data columns
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8 - clear data columns
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9 - clear data columns
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
